$wb = $excel.ActiveWorkbook

# Add a new worksheet named "metadata"
$ws = $wb.Worksheets.Add()
$ws.Name = "metadata"

# Header row (B1:G1) - bold, bordered, centered (matches the "data" sheet header style)
$hdr = $ws.Range("B1:G1")
$hdr.Font.Bold = $true
$hdr.HorizontalAlignment = -4108
$hdr.VerticalAlignment = -4160
$hdr.Borders.LineStyle = 1

$ws.Range("B1").Value = "data_name"
$ws.Range("C1").Value = "data_id"
$ws.Range("D1").Value = "data_version"
$ws.Range("E1").Value = "data_version_created"
$ws.Range("F1").Value = "panel_query_time"
$ws.Range("G1").Value = "panel_get_request"

# A2 carries the same header-style formatting as the header row
$a2 = $ws.Range("A2")
$a2.Font.Bold = $true
$a2.HorizontalAlignment = -4108
$a2.VerticalAlignment = -4160
$a2.Borders.LineStyle = 1
$a2.Value = 0

# Data row (B2:G2) - default style
$ws.Range("B2").Value = "Lung cancer pertinent cancer susceptibility"
$ws.Range("C2").Value = 218

# D2 must stay text ("1.1") instead of being auto-coerced to a number;
# force text via NumberFormat, then clear the format so no style sticks.
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "1.1"
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = "2020-05-07T14:27:40.909733Z"
$ws.Range("F2").Value = "2021-10-05 14:21:27.275007"
$ws.Range("G2").Value = "https://panelapp.genomicsengland.co.uk/api/v1/panels/218/?format=json"

# Move the new sheet after "data" so ordering matches: data, metadata
$dataSheet = $wb.Worksheets.Item("data")
$ws.Move($null, $dataSheet)
